# Weekly driver report update for 2025-04-28
# Rewrites the "Bad Drivers" and "Good Drivers" tables on the active sheet
# with the refreshed weekly figures, preserving the existing layout/styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New data (captured from this week's driver summary)
# ---------------------------------------------------------------------

# Bad drivers: Adapter-Driver, Client Count, Critical Minutes, Good Roaming Calculation (%)
$badDrivers = @(
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.0.5', 75, 7298, 91.7),
    @('Intel(R) Dual Band Wireless-AC 8260 - 20.70.18.2', 5, 1369, 97.9),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 22.100.0.3', 1, 7, 98.1),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 22.80.1.1', 1, 6, 98.1),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.0.3', 1, 5, 98.3),
    @('Intel(R) Wi-Fi 6E AX211 160MHz - 22.220.1.1', 41, 6886, 98.6),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 23.40.0.4', 2, 4, 98.7)
)
$badTotalsB = 126
$badTotalsC = 15575

# Good drivers: Adapter-Driver, Total Samples, Good Roaming Calculation (%), Driver Vintage
$goodDrivers = @(
    @('Intel(R) Dual Band Wireless-AC 8260 - 20.50.0.5', 323804, 100, $null),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3', 11128, 100, $null),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4', 486214, 99.9, '2024-11-10'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3', 18721, 99.9, '2024-07-23'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1', 69578, 99.9, '2023-08-14'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8', 338880, 99.9, '2023-05-08'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6', 143869, 99.9, '2023-01-16'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4', 287148, 99.9, '2022-11-22'),
    @('Intel(R) Dual Band Wireless-AC 8260 - 22.180.0.4', 10456, 100, '2022-10-17'),
    @('Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1', 11140, 100, '2022-08-29'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4', 96526, 99.9, '2022-08-13'),
    @('Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3', 14487, 100, '2022-05-23'),
    @('Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1', 265400, 99.9, '2022-05-01'),
    @('Intel(R) Dual Band Wireless-AC 8260 - 22.80.1.1', 123675, 100, '2021-09-11'),
    @('Intel(R) Dual Band Wireless-AC 8260 - 20.70.27.1', 18967, 100, '2021-09-11'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9', 79953, 99.9, '2021-08-18'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1', 35355, 100, '2021-04-27'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11', 67111, 100, '2021-01-19'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7', 68450, 100, '2020-10-19'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1', 15734, 99.9, '2020-09-28'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2', 65425, 100, '2020-08-05'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6', 117653, 100, '2020-01-06'),
    @('Intel(R) Dual Band Wireless-AC 8260 - 20.70.16.4', 35023, 100, '2019-12-31'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.2.1', 26241, 100, '2019-12-14'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1', 56018, 100, '2019-12-14'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2', 90508, 99.9, '2019-08-31'),
    @('Intel(R) Dual Band Wireless-AC 8260 - 20.70.12.5', 197997, 99.9, '2019-08-25'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1', 13016, 100, '2019-07-29'),
    @('Intel(R) Wi-Fi 6 AX200 160MHz - 21.10.1.2', 52515, 100, '2019-04-23'),
    @('Intel(R) Dual Band Wireless-AC 8260 - 20.70.5.2', 160536, 99.9, '2018-11-25')
)

# ---------------------------------------------------------------------
# Adjust table sizes first so everything below lands on the right rows.
# The "Bad Drivers" table shrinks from 9 data rows to 7 (remove 2 rows
# right above the old "Totals:" row); the "Good Drivers" table grows
# from 29 data rows to 30 (add 1 row right above the first data row).
# ---------------------------------------------------------------------

# Bad drivers used to occupy rows 3-11 (9 rows); now only 3-9 (7 rows).
$ws.Range("A10:A11").EntireRow.Delete()

# Good drivers used to start at row 20 (title row 18, header row 19);
# now needs one extra data row, so insert a row before the old first
# data row (which is now row 18 after the delete above).
$ws.Range("A18:A18").EntireRow.Insert()

# The inserted row doesn't automatically pick up the data-row number
# formatting used throughout the rest of the "Good Drivers" table, so
# apply it explicitly to the new row 18.
$ws.Cells.Item(18, 2).NumberFormat = "#,##0"
$ws.Cells.Item(18, 2).HorizontalAlignment = -4152
$ws.Cells.Item(18, 4).HorizontalAlignment = -4152
$ws.Cells.Item(18, 5).HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# Write the "Bad Drivers" table (title row 1 / header row 2 unchanged)
# ---------------------------------------------------------------------

$row = 3
foreach ($d in $badDrivers) {
    $ws.Cells.Item($row, 1).Value = $d[0]
    $ws.Cells.Item($row, 2).Value = $d[1]
    $ws.Cells.Item($row, 3).Value = $d[2]
    $ws.Cells.Item($row, 4).Value = $d[3]
    $row = $row + 1
}

# Totals row (row 10)
$ws.Cells.Item($row, 1).Value = "Totals:"
$ws.Cells.Item($row, 2).Value = $badTotalsB
$ws.Cells.Item($row, 3).Value = $badTotalsC

# ---------------------------------------------------------------------
# Write the "Good Drivers" table (title row 16 / header row 17 unchanged)
# ---------------------------------------------------------------------

$row = 18
foreach ($d in $goodDrivers) {
    $ws.Cells.Item($row, 1).Value = $d[0]
    $ws.Cells.Item($row, 2).Value = $d[1]
    $ws.Cells.Item($row, 4).Value = $d[2]
    if ($d[3] -eq $null) {
        $ws.Cells.Item($row, 5).ClearContents()
    } else {
        # Leading apostrophe forces text so the yyyy-mm-dd vintage stays a
        # literal string instead of being auto-converted to a date serial.
        $ws.Cells.Item($row, 5).Value = "'" + $d[3]
    }
    $row = $row + 1
}
